$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the value in C13 (was 0.286, corrected to 0.285)
$ws.Range("C13").Select() | Out-Null
$ws.Range("C13").Value = 0.285

# Pressing Enter after editing C13 moves the active selection down to C14
$ws.Range("C14").Select() | Out-Null
